$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change C9 from "00-00-00_00-15-00.mov" to the new segment value
$ws.Range("C9").Value = "00-15-00_00-21-41.mov"

# Add new row 10 with data: id=9, program_id=4, segment=00-00-00_00-15-00.mov
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "00-00-00_00-15-00.mov"

# Update selection to match target state
$ws.Range("C15").Select()
